# Update scenario cell-reference labels in the Init sheet from row 36 to row 39
# (init files updated to add more scenarios).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

$ws.Range("D5").Value  = "A39"
$ws.Range("D6").Value  = "B39"
$ws.Range("D7").Value  = "C39"
$ws.Range("D8").Value  = "G39"
$ws.Range("D9").Value  = "H39"
$ws.Range("D10").Value = "I39"
$ws.Range("D11").Value = "J39"
